$wb = $excel.ActiveWorkbook

# Activate the "Rules" worksheet (it is the one with the typo fix)
$ws = $wb.Worksheets.Item("Rules")
$ws.Activate()

# Fix the typo "retirn bar.foo;" -> "return bar.foo;"
$ws.Range("B14").Value = "return bar.foo;"

# Move the selection to the cell right below, like a user pressing Enter after editing
$ws.Range("B15").Select()

$wb.Save()
